$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the B:F (and G) values of rows 3..11 down by one, taking the value
# that used to sit in the row above (this matches a re-indexing of the
# quarterly rows). We must read all the "old" values first, since we will
# be overwriting them as we go.

$oldValues = @{}
for ($r = 2; $r -le 11; $r++) {
    $oldValues[$r] = @(
        $ws.Cells.Item($r, 2).Value2,
        $ws.Cells.Item($r, 3).Value2,
        $ws.Cells.Item($r, 4).Value2,
        $ws.Cells.Item($r, 5).Value2,
        $ws.Cells.Item($r, 6).Value2,
        $ws.Cells.Item($r, 7).Value2
    )
}

# Rows 3..11 take on the values previously held by the row directly above them.
for ($r = 11; $r -ge 3; $r--) {
    $prev = $oldValues[$r - 1]
    $ws.Cells.Item($r, 2).Value2 = $prev[0]
    $ws.Cells.Item($r, 3).Value2 = $prev[1]
    $ws.Cells.Item($r, 4).Value2 = $prev[2]
    $ws.Cells.Item($r, 5).Value2 = $prev[3]
    $ws.Cells.Item($r, 6).Value2 = $prev[4]
    $ws.Cells.Item($r, 7).Value2 = $prev[5]
}

# Row 2 (Q0) gets the freshly re-computed error metrics.
$ws.Cells.Item(2, 2).Value2 = 0.08298509261913574
$ws.Cells.Item(2, 3).Value2 = 0.2782259881569277
$ws.Cells.Item(2, 4).Value2 = 0.1147754943926806
$ws.Cells.Item(2, 5).Value2 = 0.3387853219852958
$ws.Cells.Item(2, 6).Value2 = 0.3408637854649495
$ws.Cells.Item(2, 7).Value2 = 14
